$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1808
$ws.Range("I29").Value = 734.3333
$ws.Range("J29").Value = 2344.8333
$ws.Range("K29").Value = 2202.9999
$ws.Range("L29").Value = 7034.499899999999
$ws.Range("M29").Value = -1921.9999
$ws.Range("N29").Value = -7596.499899999999

$ws.Range("H38").Value = 242.3
$ws.Range("I38").Value = 158.11111
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 474.33333
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -102.33333
$ws.Range("N38").Value = -3744

$ws.Range("H62").Value = 3382.5715
$ws.Range("I62").Value = 3935.8
$ws.Range("J62").Value = 1999.5
$ws.Range("K62").Value = 3935.8
$ws.Range("L62").Value = 1999.5
$ws.Range("M62").Value = -3311.8
$ws.Range("N62").Value = -3247.5

$ws.Range("H65").Value = 3382.5715
$ws.Range("I65").Value = 3935.8
$ws.Range("J65").Value = 1999.5
$ws.Range("K65").Value = 19679
$ws.Range("L65").Value = 9997.5
$ws.Range("M65").Value = -16559
$ws.Range("N65").Value = -16237.5

$ws.Range("H82").Value = 783
$ws.Range("I82").Value = 783
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2349
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1943
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 783
$ws.Range("I85").Value = 783
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2349
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -945
$ws.Range("N85").ClearContents()

$ws.Range("H97").Value = 2320.4707
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 2320.4707
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6961.4121
$ws.Range("N97").Value = -7953.4121

$ws.Range("H106").Value = 1771.24
$ws.Range("I106").Value = 1816.3043
$ws.Range("J106").Value = 1253
$ws.Range("K106").Value = 1816.3043
$ws.Range("L106").Value = 1253
$ws.Range("M106").Value = -1185.3043

$ws.Range("H111").Value = 2926.4285
$ws.Range("I111").Value = 1379.4
$ws.Range("J111").Value = 4332.8184
$ws.Range("K111").Value = 4138.200000000001
$ws.Range("L111").Value = 12998.4552
$ws.Range("M111").Value = -1071.200000000001

$ws.Range("H112").Value = 2346.257
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 2603.9666
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 7811.899800000001
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -10027.8998

$ws.Range("H121").Value = 1232.3334
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1232.3334
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3697.0002
$ws.Range("N121").Value = -7191.0002

$ws.Range("H126").Value = 29995.334
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 29995.334
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 29995.334
$ws.Range("N126").Value = -39875.334

$ws.Range("H128").Value = 35980.727
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 35980.727
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 35980.727
$ws.Range("N128").Value = -45940.727

$ws.Range("H133").Value = 34999.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 34999.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 34999.332
$ws.Range("N133").Value = -45119.332

$ws.Range("H138").Value = 2244.149
$ws.Range("I138").Value = 1818.3334
$ws.Range("J138").Value = 2258.1868
$ws.Range("K138").Value = 5455.0002
$ws.Range("L138").Value = 6774.5604
$ws.Range("M138").Value = -315.0002000000004
$ws.Range("N138").Value = -17054.5604

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1606.6875
$ws.Range("I2").Value = 905.7857
$ws.Range("J2").Value = 6513
$ws.Range("K2").Value = 905.7857
$ws.Range("L2").Value = 6513
$ws.Range("M2").Value = -792.7857

$ws.Range("H32").Value = 4493.4546
$ws.Range("I32").Value = 4346.129
$ws.Range("J32").Value = 6777
$ws.Range("K32").Value = 4346.129
$ws.Range("L32").Value = 6777
$ws.Range("M32").Value = -4059.129
$ws.Range("N32").Value = -7351

$ws.Range("H116").Value = 1606.6875
$ws.Range("I116").Value = 905.7857
$ws.Range("J116").Value = 6513
$ws.Range("K116").Value = 905.7857
$ws.Range("L116").Value = 6513
$ws.Range("M116").Value = 1388.2143

$ws.Range("H132").Value = 3090.7812
$ws.Range("I132").Value = 2887.261
$ws.Range("J132").Value = 3610.889
$ws.Range("K132").Value = 8661.782999999999
$ws.Range("L132").Value = 10832.667
$ws.Range("M132").Value = -6131.782999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1606.6875
$ws.Range("I3").Value = 905.7857
$ws.Range("J3").Value = 6513
$ws.Range("K3").Value = 905.7857
$ws.Range("L3").Value = 6513
$ws.Range("M3").Value = -791.7857

$ws.Range("H134").Value = 10036.917
$ws.Range("I134").Value = 1744.4
$ws.Range("J134").Value = 51499.5
$ws.Range("K134").Value = 5233.200000000001
$ws.Range("L134").Value = 154498.5
$ws.Range("M134").Value = -2698.200000000001
$ws.Range("N134").Value = -159568.5

$ws.Range("H140").Value = 23250.588
$ws.Range("I140").Value = 20780
$ws.Range("J140").Value = 23405
$ws.Range("K140").Value = 20780
$ws.Range("L140").Value = 23405
$ws.Range("M140").Value = -15600
$ws.Range("N140").Value = -33765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1135.902
$ws.Range("I31").Value = 808.0417
$ws.Range("J31").Value = 1427.3334
$ws.Range("K31").Value = 808.0417
$ws.Range("L31").Value = 1427.3334
$ws.Range("M31").Value = -513.0417
$ws.Range("N31").Value = -2017.3334

$ws.Range("H34").Value = 1135.902
$ws.Range("I34").Value = 808.0417
$ws.Range("J34").Value = 1427.3334
$ws.Range("K34").Value = 808.0417
$ws.Range("L34").Value = 1427.3334
$ws.Range("M34").Value = -606.0417
$ws.Range("N34").Value = -1831.3334

$ws.Range("H111").Value = 3351
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 3351
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 3351
$ws.Range("N111").Value = -11531

$ws.Range("H132").Value = 3125.7036
$ws.Range("I132").Value = 3070.4888
$ws.Range("J132").Value = 3401.7778
$ws.Range("K132").Value = 9211.466400000001
$ws.Range("L132").Value = 10205.3334
$ws.Range("M132").Value = -6681.466400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 897.7778
$ws.Range("I98").Value = 340.2
$ws.Range("J98").Value = 1594.75
$ws.Range("K98").Value = 1020.6
$ws.Range("L98").Value = 4784.25
$ws.Range("M98").Value = 477.4000000000001
$ws.Range("N98").Value = -7780.25

$ws.Range("H107").Value = 4209.643
$ws.Range("I107").Value = 639.4
$ws.Range("J107").Value = 8329.154
$ws.Range("K107").Value = 1918.2
$ws.Range("L107").Value = 24987.462
$ws.Range("M107").Value = 1.800000000000182
$ws.Range("N107").Value = -28827.462

$ws.Range("H113").Value = 707.5714
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 715.8461
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2147.5383
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6487.5383

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28130000
$ws.Range("I70").Value = 25004710
$ws.Range("J70").Value = 33338816
$ws.Range("K70").Value = 25004710
$ws.Range("L70").Value = 33338816
$ws.Range("M70").Value = -25004440
$ws.Range("N70").Value = -33339356

$ws.Range("H73").Value = 28130000
$ws.Range("I73").Value = 25004710
$ws.Range("J73").Value = 33338816
$ws.Range("K73").Value = 25004710
$ws.Range("L73").Value = 33338816
$ws.Range("M73").Value = -25003774
$ws.Range("N73").Value = -33340688

$ws.Range("H113").Value = 1962.9231
$ws.Range("I113").Value = 1089.75
$ws.Range("J113").Value = 3360
$ws.Range("K113").Value = 1089.75
$ws.Range("L113").Value = 3360
$ws.Range("M113").Value = 1080.25

$ws.Range("H132").Value = 2439.25
$ws.Range("I132").Value = 2130.36
$ws.Range("J132").Value = 3542.4285
$ws.Range("K132").Value = 6391.08
$ws.Range("L132").Value = 10627.2855
$ws.Range("M132").Value = -3861.08
$ws.Range("N132").Value = -15687.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10872710
$ws.Range("I122").Value = 22729536
$ws.Range("J122").Value = 3953.9167
$ws.Range("K122").Value = 68188608
$ws.Range("L122").Value = 11861.7501
$ws.Range("M122").Value = -68186158
$ws.Range("N122").Value = -16761.7501

$ws.Range("H136").Value = 1278.9166
$ws.Range("I136").Value = 1241.2632
$ws.Range("J136").Value = 1422
$ws.Range("K136").Value = 3723.7896
$ws.Range("L136").Value = 4266
$ws.Range("M136").Value = -1173.7896
$ws.Range("N136").Value = -9366

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 482.45
$ws.Range("I113").Value = 284.26666
$ws.Range("J113").Value = 1077
$ws.Range("K113").Value = 852.79998
$ws.Range("L113").Value = 3231
$ws.Range("M113").Value = 1317.20002

$ws.Range("H122").Value = 10401407
$ws.Range("I122").Value = 10834674
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 32504022
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -32501572
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 4623.174
$ws.Range("I132").Value = 5408.6665
$ws.Range("J132").Value = 3150.375
$ws.Range("K132").Value = 16225.9995
$ws.Range("L132").Value = 9451.125
$ws.Range("M132").Value = -13695.9995
$ws.Range("N132").Value = -14511.125
